# Auto commit at 2025-08-22 7:34:44.15
# Updates Metrics values, propagates them as formulas on the "today" sheet,
# and refreshes the saved selection state on the affected sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metrics sheet: refresh the 12 summary values (B2:B13)
# ---------------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 355177.95
$metrics.Range("B3").Value  = 304835.11
$metrics.Range("B4").Value  = 112218.63
$metrics.Range("B5").Value  = 13972
$metrics.Range("B6").Value  = 3750806.52
$metrics.Range("B7").Value  = 3183549.77
$metrics.Range("B8").Value  = 1075861.19
$metrics.Range("B9").Value  = 144660
$metrics.Range("B10").Value = 32216130.32
$metrics.Range("B11").Value = 19213419.84
$metrics.Range("B12").Value = 11357570.08
$metrics.Range("B13").Value = 1242287

# ---------------------------------------------------------------------------
# "today" sheet: B11:B22 now pull their numbers from Metrics via formulas
# instead of static values.
# ---------------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.Range("B11").Formula = "=Metrics!B2"
$today.Range("B12").Formula = "=Metrics!B3"
$today.Range("B13").Formula = "=Metrics!B4"
$today.Range("B14").Formula = "=Metrics!B5"
$today.Range("B15").Formula = "=Metrics!B6"
$today.Range("B16").Formula = "=Metrics!B7"
$today.Range("B17").Formula = "=Metrics!B8"

# B18 also switches its number format from "177" (#,##0.00_ ) to "4" (#,##0.00)
$today.Range("B18").NumberFormat = "#,##0.00"
$today.Range("B18").Formula = "=Metrics!B9"

$today.Range("B19").Formula = "=Metrics!B10"
$today.Range("B20").Formula = "=Metrics!B11"
$today.Range("B21").Formula = "=Metrics!B12"

# B22 also switches its number format from "177" (#,##0.00_ ) to "4" (#,##0.00)
$today.Range("B22").NumberFormat = "#,##0.00"
$today.Range("B22").Formula = "=Metrics!B13"

# E15:E22 pick up the "#,##0.00_ " number format that E11:E14 already had
$today.Range("E15:E22").NumberFormat = "#,##0.00_ "

# E11 becomes a formula referencing B11 directly
$today.Range("E11").Formula = "=B11"
# E12:E22 becomes a (shared) formula referencing the B column one row at a time
$today.Range("E12:E22").Formula = "=B12"

# F11's formula now guards against E11/B11 mismatches
$today.Range("F11").Formula = '=IF(E11=B11,E11+B3,"")'

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
# Update "today"'s selection (and implicitly drop its stale topLeftCell)
# before re-selecting on Metrics, so Metrics remains the active tab.
$today.Range("C12").Select()
$metrics.Range("C15").Select()

Write-Host "Applied Metrics/today updates"
